$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1090
$ws.Range("I40").Value = 800
$ws.Range("J40").Value = 1105.2632
$ws.Range("K40").Value = 800
$ws.Range("L40").Value = 1105.2632
$ws.Range("M40").Value = -625
$ws.Range("N40").Value = -1455.2632

$ws.Range("H76").Value = 3786.3845
$ws.Range("I76").Value = 3600.4285
$ws.Range("J76").Value = 4003.3333
$ws.Range("K76").Value = 3600.4285
$ws.Range("L76").Value = 4003.3333
$ws.Range("M76").Value = -3285.4285
$ws.Range("N76").Value = -4633.3333

$ws.Range("H79").Value = 3786.3845
$ws.Range("I79").Value = 3600.4285
$ws.Range("J79").Value = 4003.3333
$ws.Range("K79").Value = 3600.4285
$ws.Range("L79").Value = 4003.3333
$ws.Range("M79").Value = -2508.4285
$ws.Range("N79").Value = -6187.3333

$ws.Range("H137").Value = 2566741.8
$ws.Range("I137").Value = 5884698
$ws.Range("J137").Value = 2866.5454
$ws.Range("K137").Value = 17654094
$ws.Range("L137").Value = 8599.636200000001
$ws.Range("M137").Value = -17651544
$ws.Range("N137").Value = -13699.6362

$ws.Range("H138").Value = 4259805.5
$ws.Range("I138").Value = 3632.3333
$ws.Range("J138").Value = 5719065
$ws.Range("K138").Value = 10896.9999
$ws.Range("L138").Value = 17157195
$ws.Range("M138").Value = -5756.999899999999
$ws.Range("N138").Value = -17167475

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6717820.5
$ws.Range("I32").Value = 51811.06
$ws.Range("J32").Value = 20883090
$ws.Range("K32").Value = 51811.06
$ws.Range("L32").Value = 20883090
$ws.Range("M32").Value = -51524.06
$ws.Range("N32").Value = -20883664

$ws.Range("H37").Value = 13833.333
$ws.Range("J37").Value = 38000
$ws.Range("L37").Value = 38000
$ws.Range("N37").Value = -38546

$ws.Range("H61").Value = 31315060
$ws.Range("I61").Value = 41710064
$ws.Range("J61").Value = 130039.25
$ws.Range("K61").Value = 41710064
$ws.Range("L61").Value = 130039.25
$ws.Range("M61").Value = -41709852
$ws.Range("N61").Value = -130463.25

$ws.Range("H136").Value = 31315060
$ws.Range("I136").Value = 41710064
$ws.Range("J136").Value = 130039.25
$ws.Range("K136").Value = 125130192
$ws.Range("L136").Value = 390117.75
$ws.Range("M136").Value = -125127642
$ws.Range("N136").Value = -395217.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 18392.715
$ws.Range("I86").Value = 24348.4
$ws.Range("J86").Value = 3503.5
$ws.Range("K86").Value = 24348.4
$ws.Range("L86").Value = 3503.5
$ws.Range("M86").Value = -23225.4
$ws.Range("N86").Value = -5749.5

$ws.Range("H89").Value = 18392.715
$ws.Range("I89").Value = 24348.4
$ws.Range("J89").Value = 3503.5
$ws.Range("K89").Value = 121742
$ws.Range("L89").Value = 17517.5
$ws.Range("M89").Value = -116126
$ws.Range("N89").Value = -28749.5

$ws.Range("H137").Value = 50000
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 60219.832
$ws.Range("I31").Value = 34610.13
$ws.Range("J31").Value = 219000
$ws.Range("K31").Value = 34610.13
$ws.Range("L31").Value = 219000
$ws.Range("M31").Value = -34315.13
$ws.Range("N31").Value = -219590

$ws.Range("H34").Value = 60219.832
$ws.Range("I34").Value = 34610.13
$ws.Range("J34").Value = 219000
$ws.Range("K34").Value = 34610.13
$ws.Range("L34").Value = 219000
$ws.Range("M34").Value = -34408.13
$ws.Range("N34").Value = -219404

$ws.Range("H62").Value = 2966.5833
$ws.Range("I62").Value = 2480
$ws.Range("J62").Value = 3314.1428
$ws.Range("K62").Value = 2480
$ws.Range("L62").Value = 3314.1428
$ws.Range("M62").Value = -1856
$ws.Range("N62").Value = -4562.1428

$ws.Range("H65").Value = 2966.5833
$ws.Range("I65").Value = 2480
$ws.Range("J65").Value = 3314.1428
$ws.Range("K65").Value = 12400
$ws.Range("L65").Value = 16570.714
$ws.Range("M65").Value = -9280
$ws.Range("N65").Value = -22810.714

$ws.Range("H132").Value = 26852.4
$ws.Range("I132").Value = 1706.963
$ws.Range("J132").Value = 79077.53999999999
$ws.Range("K132").Value = 5120.889
$ws.Range("L132").Value = 237232.62
$ws.Range("M132").Value = -2590.889
$ws.Range("N132").Value = -242292.62

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3368885.8
$ws.Range("I4").Value = 200054
$ws.Range("J4").Value = 4002652
$ws.Range("K4").Value = 600162
$ws.Range("L4").Value = 12007956
$ws.Range("M4").Value = -600050
$ws.Range("N4").Value = -12008180

$ws.Range("H5").Value = 643.3514
$ws.Range("I5").Value = 414.96295
$ws.Range("J5").Value = 1260
$ws.Range("K5").Value = 1244.88885
$ws.Range("L5").Value = 3780
$ws.Range("M5").Value = -1132.88885
$ws.Range("N5").Value = -4004

$ws.Range("H12").Value = 117.48485
$ws.Range("I12").Value = 81.666664
$ws.Range("J12").Value = 147.33333
$ws.Range("K12").Value = 244.999992
$ws.Range("L12").Value = 441.99999
$ws.Range("M12").Value = -71.99999199999999
$ws.Range("N12").Value = -787.99999

$ws.Range("H46").Value = 2135.0908
$ws.Range("I46").Value = 271.5
$ws.Range("J46").Value = 3200
$ws.Range("K46").Value = 814.5
$ws.Range("L46").Value = 9600
$ws.Range("M46").Value = -723.5
$ws.Range("N46").Value = -9782

$ws.Range("H122").Value = 1251.4166
$ws.Range("I122").Value = 451
$ws.Range("J122").Value = 1365.762
$ws.Range("K122").Value = 4059
$ws.Range("L122").Value = 12291.858
$ws.Range("M122").Value = -1609
$ws.Range("N122").Value = -17191.858

$ws.Range("H131").Value = 1053.2839
$ws.Range("J131").Value = 1097.5466
$ws.Range("L131").Value = 3292.6398
$ws.Range("N131").Value = -13372.6398

$ws.Range("H135").Value = 643.3514
$ws.Range("I135").Value = 414.96295
$ws.Range("J135").Value = 1260
$ws.Range("K135").Value = 3734.66655
$ws.Range("L135").Value = 11340
$ws.Range("M135").Value = -1199.66655
$ws.Range("N135").Value = -16410

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 29306.324
$ws.Range("I70").Value = 47737.086
$ws.Range("J70").Value = 4370.5884
$ws.Range("K70").Value = 47737.086
$ws.Range("L70").Value = 4370.5884
$ws.Range("M70").Value = -47467.086
$ws.Range("N70").Value = -4910.5884

$ws.Range("H73").Value = 29306.324
$ws.Range("I73").Value = 47737.086
$ws.Range("J73").Value = 4370.5884
$ws.Range("K73").Value = 47737.086
$ws.Range("L73").Value = 4370.5884
$ws.Range("M73").Value = -46801.086
$ws.Range("N73").Value = -6242.5884

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 29995
$ws.Range("J106").Value = 29995
$ws.Range("L106").Value = 29995
$ws.Range("N106").Value = -32519

$ws.Range("H132").Value = 40803.73
$ws.Range("I132").Value = 2143.3914
$ws.Range("J132").Value = 337199.66
$ws.Range("K132").Value = 6430.174199999999
$ws.Range("L132").Value = 1011598.98
$ws.Range("M132").Value = -3900.174199999999
$ws.Range("N132").Value = -1016658.98
